# GIT_memo.xlsx update
# - "Удаление директории" row: command text corrected from "rm Directory" to "rm -R Directory"
# - New section "GIT. Настройка SSH подлючения" (ssh-keygen) added after the "GIT. Ветки" table
# - New section "GIT. Публикация" (git clone) added after that
#
# Constants used (standard Excel COM values, written literally since no
# type library is loaded in this host):
#   xlCenter          = -4108   (horizontal/vertical alignment)
#   xlContinuous      =  1      (Borders.LineStyle)
#   xlThin            =  2      (Borders.Weight)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Fix the existing "Удаление директории" row: the command cell used to say
#    "rm Directory" - it's corrected to "rm -R Directory".
# ---------------------------------------------------------------------------
$ws.Cells.Item(39, 2).Value = "rm -R Directory"

# ---------------------------------------------------------------------------
# Helper-ish constants
# ---------------------------------------------------------------------------
$vCenter = -4108

# ---------------------------------------------------------------------------
# 2. Blank spacer row under the "GIT. Ветки" table (row 92)
# ---------------------------------------------------------------------------
$ws.Range("A92:D92").VerticalAlignment = $vCenter

# ---------------------------------------------------------------------------
# 3. New section title: "GIT. Настройка SSH подлючения"
# ---------------------------------------------------------------------------
$titleRange = $ws.Range("A93")
$titleRange.Value = "GIT. Настройка SSH подлючения"
$titleRange.Font.Bold = $true
$titleRange.VerticalAlignment = $vCenter
$ws.Range("B93:D93").VerticalAlignment = $vCenter

# Row 94 - blank formatting-only row between the title and the note
$ws.Range("A94:D94").VerticalAlignment = $vCenter

# ---------------------------------------------------------------------------
# 4. Note line under the title
# ---------------------------------------------------------------------------
$noteRange = $ws.Range("A95")
$noteRange.Value = "Примечание: SSH подключение - настройка уникального ключа"
$noteRange.Font.Bold = $false
$noteRange.VerticalAlignment = $vCenter
$ws.Range("B95:D95").VerticalAlignment = $vCenter

# Row 96 - blank spacer row before the table
$ws.Range("A96:D96").VerticalAlignment = $vCenter

# ---------------------------------------------------------------------------
# 5. Table header (row 97) - same look as every other table on the sheet
# ---------------------------------------------------------------------------
function Set-MemoTableHeader($rowNum) {
    $hdr = $ws.Range("A" + $rowNum + ":D" + $rowNum)
    $hdr.Font.Bold = $true
    $hdr.Borders.LineStyle = 1
    $hdr.Borders.Weight = 2
    $hdr.HorizontalAlignment = $vCenter
    $hdr.VerticalAlignment = $vCenter
    $ws.Range("D" + $rowNum).WrapText = $true

    $ws.Cells.Item($rowNum, 1).Value = "Действие "
    $ws.Cells.Item($rowNum, 2).Value = "Команда"
    $ws.Cells.Item($rowNum, 3).Value = "Параметры"
    $ws.Cells.Item($rowNum, 4).Value = "Комментарий"
}

function Set-MemoDataRow($rowNum, $action, $command, $parameters, $comment) {
    $row = $ws.Range("A" + $rowNum + ":D" + $rowNum)
    $row.Borders.LineStyle = 1
    $row.Borders.Weight = 2
    $row.VerticalAlignment = $vCenter
    $ws.Range("D" + $rowNum).WrapText = $true

    $ws.Cells.Item($rowNum, 1).Value = $action
    $ws.Cells.Item($rowNum, 2).Value = $command
    $ws.Cells.Item($rowNum, 2).Font.Bold = $true
    if ($parameters) {
        $ws.Cells.Item($rowNum, 3).Value = $parameters
    }
    $ws.Cells.Item($rowNum, 4).Value = $comment
}

Set-MemoTableHeader 97

# ---------------------------------------------------------------------------
# 6. Data row (row 98) - generating an SSH key
# ---------------------------------------------------------------------------
Set-MemoDataRow 98 "Генерация ключа" "ssh-keygen" $null "Генерирует файл с уникальным ключем. Как правило в папке пользователя. Переходим на сайт GitHub в меню Пользователь - Settings - SSH and GPG keys. Жмем кнопку New SSH key. Title записываем название ключа. Например: my computer. В поле Key скопировать содержимое файла-ключа. Жмем Add SSH key. Подтверждаем пароль."
$ws.Rows.Item(98).RowHeight = 90

# ---------------------------------------------------------------------------
# 7. New section title: "GIT. Публикация" (row 100, row 99 stays blank)
# ---------------------------------------------------------------------------
$pubTitle = $ws.Range("A100")
$pubTitle.Value = "GIT. Публикация"
$pubTitle.Font.Bold = $true
$pubTitle.VerticalAlignment = $vCenter

# ---------------------------------------------------------------------------
# 8. Table header (row 102, row 101 stays blank)
# ---------------------------------------------------------------------------
Set-MemoTableHeader 102

# ---------------------------------------------------------------------------
# 9. Data row (row 103) - cloning a repository over SSH
# ---------------------------------------------------------------------------
Set-MemoDataRow 103 "Клонирование репозитория на компьютер" "git clone <ssh connection from GitHub> <Directory>" $null "Загружает проект с удаленного репозитор я по SSH подключению"
$ws.Rows.Item(103).RowHeight = 30

# ---------------------------------------------------------------------------
# 10. Update the view: scroll/selection so the new content is in view
# ---------------------------------------------------------------------------
$ws.Range("A106").Select()
